# Generate Report for Handoff
# The file "dee81500-f7f0-463a-b33c-3679c84a8e0e" transitioned from
# "In Translation" to "Ready for handoff" with a new handoff timestamp.
# This pushes its row to the bottom of the (status-sorted) block of rows
# 6-9 on every sheet, shifting the other three rows up by one.

$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink {
    param(
        $ws,
        [string]$addr,
        [string]$newValue
    )
    $ws.Range($addr).Value = $newValue
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq ('$' + $addr.Substring(0,1) + '$' + $addr.Substring(1))) {
            $h.TextToDisplay = $newValue
        }
    }
}

# ---------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn), C (de-de), D (Latest Handoff Date)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndHyperlink $wsOverview "A6" "f1f919b1-b03e-42a3-af0e-6245ec5a2699.md"
$wsOverview.Range("B6").Value = "In Translation"
$wsOverview.Range("C6").Value = "In Translation"
$wsOverview.Range("D6").Value = "2016-27-18 12:27:23"

Set-CellAndHyperlink $wsOverview "A7" "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-26-18 12:26:04"

Set-CellAndHyperlink $wsOverview "A8" "38afd42c-a8bc-4b90-9af3-da17dafedae7.md"
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Range("D8").Value = "2016-28-18 12:28:12"

Set-CellAndHyperlink $wsOverview "A9" "dee81500-f7f0-463a-b33c-3679c84a8e0e.md"
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-30-18 12:30:53"

# ---------------------------------------------------------------
# Sheet "zh-cn": A Source File Name, B File Extension, C Status,
# D Latest Handoff File, E Latest Handoff Datetime
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndHyperlink $wsZhCn "A6" "f1f919b1-b03e-42a3-af0e-6245ec5a2699.md"
$wsZhCn.Range("C6").Value = "In Translation"
Set-CellAndHyperlink $wsZhCn "D6" "f1f919b1-b03e-42a3-af0e-6245ec5a2699.88c811680106f75e284511876256a3c1ef745e9d.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-18 12:27:20"

Set-CellAndHyperlink $wsZhCn "A7" "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "D7" "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-18 12:26:00"

Set-CellAndHyperlink $wsZhCn "A8" "38afd42c-a8bc-4b90-9af3-da17dafedae7.md"
$wsZhCn.Range("C8").Value = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "D8" "38afd42c-a8bc-4b90-9af3-da17dafedae7.96b815bb8ed19869c165d55759e1fe4f7f19be3e.zh-cn.xlf"
$wsZhCn.Range("E8").Value = "2016-03-18 12:28:10"

Set-CellAndHyperlink $wsZhCn "A9" "dee81500-f7f0-463a-b33c-3679c84a8e0e.md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
Set-CellAndHyperlink $wsZhCn "D9" "dee81500-f7f0-463a-b33c-3679c84a8e0e.14d86bff35e19cab99635e18fef74dfa734a9d1f.zh-cn.xlf"
$wsZhCn.Range("E9").Value = "2016-03-18 12:30:50"

# ---------------------------------------------------------------
# Sheet "de-de": A Source File Name, B File Extension, C Status,
# D Latest Handoff File, E Latest Handoff Datetime
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndHyperlink $wsDeDe "A6" "f1f919b1-b03e-42a3-af0e-6245ec5a2699.md"
$wsDeDe.Range("C6").Value = "In Translation"
Set-CellAndHyperlink $wsDeDe "D6" "f1f919b1-b03e-42a3-af0e-6245ec5a2699.88c811680106f75e284511876256a3c1ef745e9d.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-18 12:27:23"

Set-CellAndHyperlink $wsDeDe "A7" "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "D7" "03ca7e5d-cec3-4d5c-8b72-f3b0f7777b0a.81200effe507a49bc7034878fb2a2f18ca8e9f06.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-18 12:26:04"

Set-CellAndHyperlink $wsDeDe "A8" "38afd42c-a8bc-4b90-9af3-da17dafedae7.md"
$wsDeDe.Range("C8").Value = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "D8" "38afd42c-a8bc-4b90-9af3-da17dafedae7.96b815bb8ed19869c165d55759e1fe4f7f19be3e.de-de.xlf"
$wsDeDe.Range("E8").Value = "2016-03-18 12:28:12"

Set-CellAndHyperlink $wsDeDe "A9" "dee81500-f7f0-463a-b33c-3679c84a8e0e.md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
Set-CellAndHyperlink $wsDeDe "D9" "dee81500-f7f0-463a-b33c-3679c84a8e0e.14d86bff35e19cab99635e18fef74dfa734a9d1f.de-de.xlf"
$wsDeDe.Range("E9").Value = "2016-03-18 12:30:53"

Write-Output "Done"
